$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '66.964.35'
Set-TextValue 'E2' '  -1.02%  '
Set-TextValue 'D3' '3.519.44'
Set-TextValue 'E3' '  +0.44%  '
Set-TextValue 'E4' '  -0.06%  '
Set-TextValue 'D5' '606.75'
Set-TextValue 'E5' '  +0.04%  '
Set-TextValue 'D6' '147.49'
Set-TextValue 'E6' '  -2.88%  '
Set-TextValue 'D7' '3.518.54'
Set-TextValue 'E7' '  +0.47%  '
Set-TextValue 'E8' '  -0.05%  '
Set-TextValue 'E9' '  -2.03%  '
Set-TextValue 'E10' '  -0.61%  '
Set-TextValue 'E11' '  +3.17%  '
Set-TextValue 'E12' '  -2.03%  '
Set-TextValue 'E13' '  -0.49%  '
Set-TextValue 'D14' '4.111.63'
Set-TextValue 'E14' '  +0.39%  '
Set-TextValue 'D15' '31.62'
Set-TextValue 'E15' '  -2.13%  '
Set-TextValue 'D16' '3.517.88'
Set-TextValue 'E16' '  +0.27%  '
Set-TextValue 'D17' '66.971.44'
Set-TextValue 'E18' '  -0.03%  '
Set-TextValue 'E19' '  +8.11%  '
Set-TextValue 'D20' '6.38'
Set-TextValue 'E20' '  -1.88%  '
Set-TextValue 'D21' '15.32'
Set-TextValue 'E21' '  -1.36%  '
Set-TextValue 'D22' '436.77'
Set-TextValue 'E22' '  -2.22%  '
Set-TextValue 'D23' '0.609'
Set-TextValue 'E23' '  -3.15%  '
Set-TextValue 'D24' '79.68'
Set-TextValue 'E24' '  +1.89%  '
Set-TextValue 'D25' '3.650.58'
Set-TextValue 'E25' '  +0.20%  '
Set-TextValue 'E26' '  -0.05%  '
Set-TextValue 'D27' '0.0000121'
Set-TextValue 'E27' '  -3.60%  '
Set-TextValue 'D28' '9.81'
Set-TextValue 'E28' '  -2.68%  '
Set-TextValue 'D29' '8.33'
Set-TextValue 'E29' '  -4.96%  '
Set-TextValue 'D30' '2.51'
Set-TextValue 'E30' '  -0.04%  '
Set-TextValue 'E31' '  -3.07%  '
Set-TextValue 'E32' '  -1.33%  '
Set-TextValue 'E33' '  +1.01%  '
Set-TextValue 'D34' '25.40'
Set-TextValue 'E34' '  -0.85%  '
Set-TextValue 'D35' '3.514.13'
Set-TextValue 'E35' '  +0.58%  '
Set-TextValue 'D36' '5.94'
Set-TextValue 'E36' '  -3.25%  '
Set-TextValue 'E37' '  -3.21%  '
Set-TextValue 'E38' '  +1.14%  '
Set-TextValue 'E39' '  +0.01%  '
Set-TextValue 'E40' '  -0.07%  '
Set-TextValue 'D41' '0.0892'
Set-TextValue 'E41' '  +0.07%  '
Set-TextValue 'D42' '169.97'
Set-TextValue 'E42' '  -4.83%  '
Set-TextValue 'D43' '5.45'
Set-TextValue 'E43' '  +0.01%  '
Set-TextValue 'E44' '  -10.05%  '
Set-TextValue 'D45' '0.896'
Set-TextValue 'E45' '  +0.75%  '
Set-TextValue 'D46' '45.98'
Set-TextValue 'E46' '  -0.77%  '
Set-TextValue 'D47' '1.34'
Set-TextValue 'E47' '  +2.97%  '
Set-TextValue 'D48' '28.20'
Set-TextValue 'E48' '  -6.82%  '
Set-TextValue 'D49' '7.48'
Set-TextValue 'E49' '  -1.80%  '
Set-TextValue 'E50' '  -4.89%  '
Set-TextValue 'D51' '0.990'
Set-TextValue 'E51' '  +0.16%  '
